$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 4 RR vs PBKS - fill in the scores for rows 11-13 (matches 2, 3, 4)
$data = @{
    11 = @{ E = 20;  H = 80; K = 40;  N = 0;   Q = 100; T = 60 }
    12 = @{ E = 60;  H = 20; K = 80;  N = 100; Q = 40;  T = 0 }
    13 = @{ E = 100; H = 60; K = 0;   N = 20;  Q = 80;  T = 40 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

$excel.Calculate()
